$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 2431
$ws.Range("B2").Value = "Antônio Barros"
$ws.Range("C2").Value = "P&D"
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 45095
$ws.Range("G2").Value = 8719.27

# Row 3
$ws.Range("A3").Value = 77644
$ws.Range("B3").Value = "Sr. Eduardo da Paz"
$ws.Range("C3").Value = "TI"
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 45101
$ws.Range("G3").Value = 4302.18

# Row 4
$ws.Range("A4").Value = 58502
$ws.Range("B4").Value = "Evelyn Jesus"
$ws.Range("C4").Value = "Financeiro"
$ws.Range("D4").Value = "Consulta médica"
$ws.Range("F4").Value = 45101
$ws.Range("G4").Value = 7512.2

# Row 5
$ws.Range("A5").Value = 58997
$ws.Range("B5").Value = "Sophie das Neves"
$ws.Range("C5").Value = "Recursos Humanos"
$ws.Range("D5").Value = "Outros"
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 9952.440000000001

# Row 6
$ws.Range("A6").Value = 26528
$ws.Range("B6").Value = "Ana Beatriz Novaes"
$ws.Range("C6").Value = "P&D"
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 45106
$ws.Range("G6").Value = 8481.870000000001

# Row 7
$ws.Range("A7").Value = 57087
$ws.Range("B7").Value = "Maysa Cavalcanti"
$ws.Range("C7").Value = "P&D"
$ws.Range("D7").Value = "Consulta médica"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 45082
$ws.Range("G7").Value = 3041.16

# Row 8
$ws.Range("A8").Value = 79342
$ws.Range("B8").Value = "Eduarda da Paz"
$ws.Range("C8").Value = "Recursos Humanos"
$ws.Range("D8").Value = "Viagem de negócios"
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 45084
$ws.Range("G8").Value = 6650.43

# Row 9
$ws.Range("A9").Value = 19808
$ws.Range("B9").Value = "Luiz Henrique Santos"
$ws.Range("C9").Value = "Financeiro"
$ws.Range("D9").Value = "Doença"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 45088
$ws.Range("G9").Value = 12426.17

# Row 10
$ws.Range("A10").Value = 25161
$ws.Range("B10").Value = "Nicolas Barbosa"
$ws.Range("D10").Value = "Consulta médica"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 45106
$ws.Range("G10").Value = 8139.02

# Row 11
$ws.Range("A11").Value = 11039
$ws.Range("B11").Value = "João Pedro Silveira"
$ws.Range("C11").Value = "Engenharia"
$ws.Range("D11").Value = "Problemas pessoais"
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 45087
$ws.Range("G11").Value = 6926.14
